$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.071.28"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.760.63"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'579.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'158.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.19%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").Value = "'5.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -14.09%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "3.248.21"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'27.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "63.736.09"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "2.762.23"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'12.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "'360.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'6.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'0.551"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.88%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'65.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").Value = "'8.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "0.0₃0932"
$ws.Range("E28").Value = "  +3.77%  "
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "'1.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").Value = "'166.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("D33").Value = "'20.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("E34").Value = "  +3.78%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'6.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.93%  "
$ws.Range("D40").Value = "'4.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'330.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.17%  "
$ws.Range("D42").Value = "'39.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "'21.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").Value = "'21.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").Value = "'0.0257"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "'136.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  +0.69%  "
